# Adds the radio-station logo URL (shared string already used by G26) together
# with an external hyperlink to that same URL, in every row of column G that
# was still blank in the original schedule.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logoUrl = "https://static.wixstatic.com/media/c39cce_e474bc38cd9d46a5b7af4e985ee1892c~mv2.png/v1/fill/w_162,h_171,al_c,q_85,usm_0.66_1.00_0.01,enc_auto/logo%20radio%20sol%20-%20color%201.png"

# Row order matches the order in which the hyperlinks were (re-)created in the
# source workbook, which drives the r:id numbering (rId126 .. rId147).
$rows = @(24, 19, 23, 31, 33, 34, 35, 38, 40, 41, 29, 60, 66, 67, 69, 71, 73, 74, 76, 77, 79, 80)

foreach ($r in $rows) {
    $gCell = $ws.Range("G" + $r)
    $hCell = $ws.Range("H" + $r)

    # Put the literal URL text into the cell first (this is what the shared
    # string table records as the cell's display text).
    $gCell.Value = $logoUrl

    # Attach the external hyperlink pointing at the same URL.
    $ws.Hyperlinks.Add($gCell, $logoUrl)

    # Adding a hyperlink re-styles the cell with a brand new "Hyperlink" style;
    # restore the original style (shared with column H in the same row) so the
    # cell keeps using the workbook's existing hyperlink-style index.
    $gCell.Style = $hCell.Style
}

# Reflect the reviewer's last on-screen position: scrolled down so row 64 is
# at the top, with G80 as the active/selected cell.
$ws.Range("G80").Select()
